$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to remain text so values like "1.000" or
# "28.613.57" are not reinterpreted as numbers/dates by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.613.57"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "1.803.90"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "317.42"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "0.5550"
$ws.Range("E7").Value = "  -2.31%  "
$ws.Range("D8").Value = "0.3788"
$ws.Range("E8").Value = "  -2.27%  "
$ws.Range("D9").Value = "0.07513"
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("D10").Value = "42.33"
$ws.Range("E10").Value = "  -1.89%  "
$ws.Range("D11").Value = "1.116"
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("D12").Value = "1.0000"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").Value = "20.69"
$ws.Range("E13").Value = "  -2.21%  "
$ws.Range("D14").Value = "6.163"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").Value = "7.382"
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("D16").Value = "1.804.84"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "90.31"
$ws.Range("E17").Value = "  -1.68%  "
$ws.Range("D18").Value = "0.00001068"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "17.29"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "5.934"
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("D23").Value = "28.622.92"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("D25").Value = "2.090"
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("D26").Value = "159.04"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("D27").Value = "20.47"
$ws.Range("E27").Value = "  -1.77%  "
$ws.Range("D28").Value = "2.010.51"
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("D29").Value = "2.354"
$ws.Range("E29").Value = "  -3.13%  "
$ws.Range("D30").Value = "123.16"
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("D31").Value = "1.107"
$ws.Range("E31").Value = "  -4.14%  "
$ws.Range("D32").Value = "0.1067"
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("D33").Value = "5.653"
$ws.Range("E33").Value = "  -2.05%  "
$ws.Range("D34").Value = "3.679"
$ws.Range("E34").Value = "  +0.93%  "
$ws.Range("D35").Value = "0.06513"
$ws.Range("E35").Value = "  +6.65%  "
$ws.Range("D36").Value = "0.2253"
$ws.Range("E36").Value = "  +4.04%  "
$ws.Range("D37").Value = "0.02308"
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").Value = "8.765"
$ws.Range("E38").Value = "  -1.80%  "
$ws.Range("D39").Value = "5.026"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").Value = "11.28"
$ws.Range("E40").Value = "  -3.17%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.6244"
$ws.Range("E41").Value = "  -2.39%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "1.205"
$ws.Range("E42").Value = "  +3.99%  "
$ws.Range("D43").Value = "1.430"
$ws.Range("E43").Value = "  +3.70%  "
$ws.Range("D44").Value = "0.9998"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "13.21"
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.5856"
$ws.Range("E46").Value = "  -2.26%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.693"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").Value = "126.63"
$ws.Range("E48").Value = "  +3.95%  "
$ws.Range("D49").Value = "1.945"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").Value = "1.159"
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("D51").Value = "0.06887"
$ws.Range("E51").Value = "  +0.34%  "
